$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Förändrad" (last changed) date for every data row
# (rows 2 through 302). The whole column was bumped by one day,
# from 45178 (2023-09-09) to 45179 (2023-09-10).
$rng = $ws.Range("C2:C302")
$rng.Value2 = 45179
